# feat: add 2022-Q1 data
#
# 1. Insert a brand-new "2022-Q1" sheet (fund-holding detail, same shape as
#    the other quarterly sheets) right before the "总计" (summary) sheet.
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q1" detail sheet, cloned from "2021-Q4" so it keeps
# the same column styling (header row style, index-column style, etc.)
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $template)

$q1 = $wb.Worksheets.Item("2021-Q4 (2)")
$q1.Name = "2022-Q1"

# Drop every data row but the first, keeping just the header + row 2.
$q1.Range("A3:H9").Clear()

# Row 2: a single fund entry. Numeric-looking text values keep a leading
# apostrophe so Excel stores them as text (matching the source columns),
# not auto-converted numbers.
$q1.Range("B2").Value = "'010714"
$q1.Range("C2").Value = "东方红远见价值混合"
$q1.Range("D2").Value = "'19.70"
$q1.Range("E2").Value = "'86.34"
$q1.Range("F2").Value = "'2.93"
$q1.Range("G2").Value = "'0.5772"
$q1.Range("H2").Value = 9

# ---------------------------------------------------------------------
# Step 2: prepend a new row to the "总计" summary sheet.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Column A carries the index-column style (matches A3:A6); bring that
# back for the new A2 (Insert leaves it unset / default).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

# Insert() left the old row-2 formatting on B2:D2 (from the row that got
# pushed to row 3) — clear it so the new data row is plain, like the rest.
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.58

# Renumber the leading index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
